$wb = $excel.ActiveWorkbook

# New week label shared across all sheets
$newWeek = "14/02/2022 - 20/02/2022"

# Data to append as row 16 on each sheet: Personale scolastico, Alunni, Totale
$rows = @{
    "Infanzia"    = @(7, 170, 177)
    "Primaria"    = @(29, 315, 344)
    "Media"       = @(9, 148, 157)
    "Superiore"   = @(5, 174, 179)
    "Totale casi" = @(50, 807, 857)
}

foreach ($name in @("Infanzia", "Primaria", "Media", "Superiore", "Totale casi")) {
    $ws = $wb.Worksheets.Item($name)
    $vals = $rows[$name]

    $ws.Range("A16").Value = $newWeek
    $ws.Range("B16").Value = $vals[0]
    $ws.Range("C16").Value = $vals[1]
    $ws.Range("D16").Value = $vals[2]
}

# Update selections to match the saved workbook state
$wb.Worksheets.Item("Infanzia").Range("A16").Select()
$wb.Worksheets.Item("Primaria").Range("E16").Select()
$wb.Worksheets.Item("Media").Range("E16").Select()
$wb.Worksheets.Item("Superiore").Range("E16").Select()
$wb.Worksheets.Item("Totale casi").Range("B17").Select()

$wb.Worksheets.Item("Totale casi").Activate()
